$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.205.16"
$ws.Range("E2").Value = "  -0.38%  "
$ws.Range("D3").Value = "2.616.00"
$ws.Range("E3").Value = "  -2.73%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "'605.45"
$ws.Range("E5").Value = "  +0.14%  "
$ws.Range("D6").Value = "'145.55"
$ws.Range("E6").Value = "  +1.19%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("E8").Value = "  -0.50%  "
$ws.Range("D9").Value = "2.612.75"
$ws.Range("E9").Value = "  -2.83%  "
$ws.Range("D10").Value = "'0.107"
$ws.Range("E10").Value = "  +0.31%  "
$ws.Range("B11").Value = "Cardano"
$ws.Range("C11").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D11").Value = "'0.376"
$ws.Range("E11").Value = "  +5.52%  "
$ws.Range("B12").Value = "Toncoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D12").Value = "'5.50"
$ws.Range("E12").Value = "  -3.24%  "
$ws.Range("E13").Value = "  -0.58%  "
$ws.Range("D14").Value = "'27.19"
$ws.Range("E14").Value = "  -0.61%  "
$ws.Range("D15").Value = "3.085.39"
$ws.Range("E15").Value = "  -2.73%  "
$ws.Range("D16").Value = "63.047.27"
$ws.Range("E16").Value = "  -0.49%  "
$ws.Range("D17").Value = "'0.0000145"
$ws.Range("D18").Value = "2.597.04"
$ws.Range("E18").Value = "  -4.18%  "
$ws.Range("D19").Value = "'11.49"
$ws.Range("E19").Value = "  +0.10%  "
$ws.Range("D20").Value = "'4.52"
$ws.Range("E20").Value = "  +2.74%  "
$ws.Range("D21").Value = "'342.68"
$ws.Range("E21").Value = "  +0.95%  "
$ws.Range("D22").Value = "'6.88"
$ws.Range("E22").Value = "  -0.01%  "
$ws.Range("E23").Value = "  -0.31%  "
$ws.Range("D24").Value = "'5.72"
$ws.Range("E24").Value = "  -1.64%  "
$ws.Range("D25").Value = "'65.95"
$ws.Range("E25").Value = "  -2.51%  "
$ws.Range("E26").Value = "  +2.23%  "
$ws.Range("B27").Value = "InternetComputer(DFINITY)"
$ws.Range("C27").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D27").Value = "'9.03"
$ws.Range("E27").Value = "  +5.66%  "
$ws.Range("B28").Value = "SuiNetwork"
$ws.Range("C28").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D28").Value = "'1.59"
$ws.Range("E28").Value = "  +3.06%  "
$ws.Range("D29").Value = "'547.61"
$ws.Range("E29").Value = "  +0.64%  "
$ws.Range("D30").Value = "'0.162"
$ws.Range("E30").Value = "  -2.06%  "
$ws.Range("D31").Value = "'1.00"
$ws.Range("E31").Value = "  +0.06%  "
$ws.Range("D32").Value = "'7.90"
$ws.Range("E32").Value = "  +0.00%  "
$ws.Range("E33").Value = "  -0.37%  "
$ws.Range("D34").Value = "0.0₃0840"
$ws.Range("E34").Value = "  +3.23%  "
$ws.Range("D35").Value = "'1.75"
$ws.Range("E35").Value = "  -3.55%  "
$ws.Range("D36").Value = "'5.22"
$ws.Range("E36").Value = "  +2.57%  "
$ws.Range("D37").Value = "'168.59"
$ws.Range("E37").Value = "  -2.62%  "
$ws.Range("E38").Value = "  -0.12%  "
$ws.Range("D39").Value = "'0.401"
$ws.Range("E39").Value = "  -1.35%  "
$ws.Range("D40").Value = "'1.93"
$ws.Range("E40").Value = "  +5.87%  "
$ws.Range("D41").Value = "'18.93"
$ws.Range("E41").Value = "  -1.74%  "
$ws.Range("E42").Value = "  +0.04%  "
$ws.Range("D43").Value = "'165.49"
$ws.Range("E43").Value = "  -5.99%  "
$ws.Range("D44").Value = "'39.68"
$ws.Range("E44").Value = "  -1.04%  "
$ws.Range("D45").Value = "'3.75"
$ws.Range("E45").Value = "  -0.30%  "
$ws.Range("D46").Value = "'21.84"
$ws.Range("E46").Value = "  -1.55%  "
$ws.Range("E47").Value = "  -0.46%  "
$ws.Range("D48").Value = "'0.621"
$ws.Range("E48").Value = "  -2.65%  "
$ws.Range("E49").Value = "  +0.35%  "
$ws.Range("B50").Value = "dogwifhat"
$ws.Range("C50").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D50").Value = "'1.93"
$ws.Range("E50").Value = "  +11.75%  "
$ws.Range("B51").Value = "Stellar"
$ws.Range("C51").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D51").Value = "'0.0954"
$ws.Range("E51").Value = "  -0.90%  "
